# Auto-generated script to update odds values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


# Row 2
$ws.Cells.Item(2, 12).Value = 4.55
$ws.Cells.Item(2, 20).Value = 2.22
$ws.Cells.Item(2, 23).Value = 6.3
$ws.Cells.Item(2, 24).Value = 9.5
$ws.Cells.Item(2, 27).Value = 18.5
$ws.Cells.Item(2, 34).Value = 10
$ws.Cells.Item(2, 37).Value = 70
$ws.Cells.Item(2, 38).Value = 40
$ws.Cells.Item(2, 39).Value = 45

# Row 3
$ws.Cells.Item(3, 7).Value = 1.8
$ws.Cells.Item(3, 8).Value = 3.5
$ws.Cells.Item(3, 9).Value = 4.05
$ws.Cells.Item(3, 10).Value = 2.3
$ws.Cells.Item(3, 11).Value = 2.18
$ws.Cells.Item(3, 12).Value = 4.35
$ws.Cells.Item(3, 13).Value = 1.33
$ws.Cells.Item(3, 14).Value = 2.8
$ws.Cells.Item(3, 17).Value = 3.2
$ws.Cells.Item(3, 21).Value = 1.88
$ws.Cells.Item(3, 22).Value = 1.72
$ws.Cells.Item(3, 23).Value = 6.3
$ws.Cells.Item(3, 24).Value = 7.8
$ws.Cells.Item(3, 25).Value = 8.5
$ws.Cells.Item(3, 26).Value = 14
$ws.Cells.Item(3, 28).Value = 32
$ws.Cells.Item(3, 29).Value = 9
$ws.Cells.Item(3, 30).Value = 6.8
$ws.Cells.Item(3, 31).Value = 17.5
$ws.Cells.Item(3, 35).Value = 21
$ws.Cells.Item(3, 37).Value = 60
$ws.Cells.Item(3, 38).Value = 40
$ws.Cells.Item(3, 39).Value = 55

# Row 4
$ws.Cells.Item(4, 8).Value = 5.5
$ws.Cells.Item(4, 9).Value = 11.75
$ws.Cells.Item(4, 10).Value = 1.62
$ws.Cells.Item(4, 11).Value = 2.55
$ws.Cells.Item(4, 21).Value = 2.1
$ws.Cells.Item(4, 26).Value = 6.9
$ws.Cells.Item(4, 29).Value = 14
$ws.Cells.Item(4, 30).Value = 11.75
$ws.Cells.Item(4, 31).Value = 28
$ws.Cells.Item(4, 35).Value = 100
$ws.Cells.Item(4, 37).Value = 450

# Row 5
$ws.Cells.Item(5, 7).Value = 2.47
$ws.Cells.Item(5, 8).Value = 2.85
$ws.Cells.Item(5, 10).Value = 3.15
$ws.Cells.Item(5, 12).Value = 3.6
$ws.Cells.Item(5, 13).Value = 1.47
$ws.Cells.Item(5, 14).Value = 2.32
$ws.Cells.Item(5, 15).Value = 2.37
$ws.Cells.Item(5, 16).Value = 1.45
$ws.Cells.Item(5, 17).Value = 4.05
$ws.Cells.Item(5, 18).Value = 1.16
$ws.Cells.Item(5, 19).Value = 1.5
$ws.Cells.Item(5, 20).Value = 2.25
$ws.Cells.Item(5, 21).Value = 1.98
$ws.Cells.Item(5, 22).Value = 1.65
$ws.Cells.Item(5, 24).Value = 10.75
$ws.Cells.Item(5, 26).Value = 27
$ws.Cells.Item(5, 27).Value = 26
$ws.Cells.Item(5, 30).Value = 5.7
$ws.Cells.Item(5, 31).Value = 17
$ws.Cells.Item(5, 32).Value = 110
$ws.Cells.Item(5, 34).Value = 7.2
$ws.Cells.Item(5, 36).Value = 11.25
$ws.Cells.Item(5, 39).Value = 45

# Row 6
$ws.Cells.Item(6, 7).Value = 2.37
$ws.Cells.Item(6, 9).Value = 2.77
$ws.Cells.Item(6, 10).Value = 2.95
$ws.Cells.Item(6, 12).Value = 3.35
$ws.Cells.Item(6, 16).Value = 1.62
$ws.Cells.Item(6, 23).Value = 7.1
$ws.Cells.Item(6, 24).Value = 10.75
$ws.Cells.Item(6, 25).Value = 9.5
$ws.Cells.Item(6, 26).Value = 24
$ws.Cells.Item(6, 27).Value = 21
$ws.Cells.Item(6, 34).Value = 7.9
$ws.Cells.Item(6, 35).Value = 13
$ws.Cells.Item(6, 36).Value = 10.5
$ws.Cells.Item(6, 37).Value = 32
$ws.Cells.Item(6, 38).Value = 26
$ws.Cells.Item(6, 39).Value = 40

# Row 7
$ws.Cells.Item(7, 7).Value = 1.9
$ws.Cells.Item(7, 8).Value = 3.4
$ws.Cells.Item(7, 9).Value = 3.7
$ws.Cells.Item(7, 10).Value = 2.5
$ws.Cells.Item(7, 11).Value = 2.07
$ws.Cells.Item(7, 12).Value = 4.15
$ws.Cells.Item(7, 13).Value = 1.28
$ws.Cells.Item(7, 14).Value = 3.05
$ws.Cells.Item(7, 15).Value = 1.82
$ws.Cells.Item(7, 16).Value = 1.78
$ws.Cells.Item(7, 17).Value = 2.87
$ws.Cells.Item(7, 18).Value = 1.31
$ws.Cells.Item(7, 19).Value = 1.4
$ws.Cells.Item(7, 20).Value = 2.52
$ws.Cells.Item(7, 21).Value = 1.7
$ws.Cells.Item(7, 22).Value = 1.91
$ws.Cells.Item(7, 23).Value = 7.4
$ws.Cells.Item(7, 24).Value = 9.25
$ws.Cells.Item(7, 26).Value = 16.5
$ws.Cells.Item(7, 27).Value = 15
$ws.Cells.Item(7, 28).Value = 25
$ws.Cells.Item(7, 29).Value = 10
$ws.Cells.Item(7, 30).Value = 6.6
$ws.Cells.Item(7, 31).Value = 14.5
$ws.Cells.Item(7, 32).Value = 65
$ws.Cells.Item(7, 33).Value = 500
$ws.Cells.Item(7, 34).Value = 10.75
$ws.Cells.Item(7, 35).Value = 20
$ws.Cells.Item(7, 36).Value = 12.5

# Row 9
$ws.Cells.Item(9, 7).Value = 2.62
$ws.Cells.Item(9, 9).Value = 2.42
$ws.Cells.Item(9, 10).Value = 3.15
$ws.Cells.Item(9, 12).Value = 2.95
$ws.Cells.Item(9, 17).Value = 2.65
$ws.Cells.Item(9, 23).Value = 9.5
$ws.Cells.Item(9, 24).Value = 14
$ws.Cells.Item(9, 25).Value = 9.75
$ws.Cells.Item(9, 26).Value = 29
$ws.Cells.Item(9, 27).Value = 21
$ws.Cells.Item(9, 28).Value = 27
$ws.Cells.Item(9, 35).Value = 13
$ws.Cells.Item(9, 36).Value = 9.25
$ws.Cells.Item(9, 37).Value = 26
$ws.Cells.Item(9, 38).Value = 18.5
$ws.Cells.Item(9, 39).Value = 25

# Row 10
$ws.Cells.Item(10, 7).Value = 1.95
$ws.Cells.Item(10, 8).Value = 3.4
$ws.Cells.Item(10, 9).Value = 4
$ws.Cells.Item(10, 10).Value = 2.63
$ws.Cells.Item(10, 12).Value = 4.33
$ws.Cells.Item(10, 21).Value = 1.8
$ws.Cells.Item(10, 22).Value = 1.91
$ws.Cells.Item(10, 23).Value = 7
$ws.Cells.Item(10, 24).Value = 9
$ws.Cells.Item(10, 25).Value = 9
$ws.Cells.Item(10, 26).Value = 17
$ws.Cells.Item(10, 27).Value = 17
$ws.Cells.Item(10, 30).Value = 6.5
$ws.Cells.Item(10, 31).Value = 15
$ws.Cells.Item(10, 32).Value = 51
$ws.Cells.Item(10, 33).Value = 251
$ws.Cells.Item(10, 34).Value = 11
$ws.Cells.Item(10, 35).Value = 21
$ws.Cells.Item(10, 36).Value = 13
$ws.Cells.Item(10, 37).Value = 41
$ws.Cells.Item(10, 38).Value = 34
$ws.Cells.Item(10, 39).Value = 41

# Row 11
$ws.Cells.Item(11, 7).Value = 2.2
$ws.Cells.Item(11, 9).Value = 3.4
$ws.Cells.Item(11, 10).Value = 3
$ws.Cells.Item(11, 12).Value = 4
$ws.Cells.Item(11, 21).Value = 2.1
$ws.Cells.Item(11, 22).Value = 1.67
$ws.Cells.Item(11, 24).Value = 9.5
$ws.Cells.Item(11, 26).Value = 21
$ws.Cells.Item(11, 33).Value = 1250
$ws.Cells.Item(11, 34).Value = 8
$ws.Cells.Item(11, 35).Value = 15
$ws.Cells.Item(11, 36).Value = 13

# Row 12
$ws.Cells.Item(12, 8).Value = 2.9
$ws.Cells.Item(12, 9).Value = 2.5
$ws.Cells.Item(12, 42).Value = 1.85
$ws.Cells.Item(12, 43).Value = 2

# Row 13
$ws.Cells.Item(13, 7).Value = 2.87
$ws.Cells.Item(13, 9).Value = 2.3
$ws.Cells.Item(13, 10).Value = 3.5
$ws.Cells.Item(13, 12).Value = 2.92
$ws.Cells.Item(13, 13).Value = 1.29
$ws.Cells.Item(13, 14).Value = 3.25
$ws.Cells.Item(13, 15).Value = 1.87
$ws.Cells.Item(13, 16).Value = 1.83
$ws.Cells.Item(13, 17).Value = 3.05
$ws.Cells.Item(13, 18).Value = 1.33
$ws.Cells.Item(13, 21).Value = 1.7
$ws.Cells.Item(13, 22).Value = 2.05
$ws.Cells.Item(13, 23).Value = 9.75
$ws.Cells.Item(13, 24).Value = 16
$ws.Cells.Item(13, 25).Value = 10.25
$ws.Cells.Item(13, 26).Value = 37
$ws.Cells.Item(13, 27).Value = 24
$ws.Cells.Item(13, 28).Value = 30
$ws.Cells.Item(13, 29).Value = 7.2
$ws.Cells.Item(13, 31).Value = 13
$ws.Cells.Item(13, 32).Value = 55
$ws.Cells.Item(13, 33).Value = 400
$ws.Cells.Item(13, 34).Value = 8.25
$ws.Cells.Item(13, 35).Value = 11.75
$ws.Cells.Item(13, 36).Value = 9
$ws.Cells.Item(13, 37).Value = 24
$ws.Cells.Item(13, 38).Value = 18.5
$ws.Cells.Item(13, 41).Value = 7.2

# Row 14
$ws.Cells.Item(14, 7).Value = 2.82
$ws.Cells.Item(14, 8).Value = 3.45
$ws.Cells.Item(14, 9).Value = 2.18
$ws.Cells.Item(14, 10).Value = 3.3
$ws.Cells.Item(14, 12).Value = 2.75
$ws.Cells.Item(14, 23).Value = 11.5
$ws.Cells.Item(14, 24).Value = 17
$ws.Cells.Item(14, 25).Value = 10.25
$ws.Cells.Item(14, 26).Value = 35
$ws.Cells.Item(14, 27).Value = 22
$ws.Cells.Item(14, 31).Value = 12
$ws.Cells.Item(14, 34).Value = 9.5
$ws.Cells.Item(14, 35).Value = 12
$ws.Cells.Item(14, 37).Value = 22

# Row 15
$ws.Cells.Item(15, 7).Value = 2.3
$ws.Cells.Item(15, 8).Value = 3.15
$ws.Cells.Item(15, 9).Value = 2.9
$ws.Cells.Item(15, 10).Value = 2.95
$ws.Cells.Item(15, 11).Value = 2.07
$ws.Cells.Item(15, 12).Value = 3.5
$ws.Cells.Item(15, 13).Value = 1.29
$ws.Cells.Item(15, 14).Value = 3.3
$ws.Cells.Item(15, 15).Value = 1.85
$ws.Cells.Item(15, 16).Value = 1.85
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 18).Value = 1.34
$ws.Cells.Item(15, 19).Value = 1.42
$ws.Cells.Item(15, 20).Value = 2.67
$ws.Cells.Item(15, 21).Value = 1.65
$ws.Cells.Item(15, 22).Value = 2.1
$ws.Cells.Item(15, 23).Value = 8.25
$ws.Cells.Item(15, 24).Value = 11.75
$ws.Cells.Item(15, 26).Value = 24
$ws.Cells.Item(15, 27).Value = 18.5
$ws.Cells.Item(15, 29).Value = 7.3
$ws.Cells.Item(15, 30).Value = 6.2
$ws.Cells.Item(15, 31).Value = 12.5
$ws.Cells.Item(15, 32).Value = 55
$ws.Cells.Item(15, 34).Value = 10
$ws.Cells.Item(15, 35).Value = 16
$ws.Cells.Item(15, 36).Value = 10.25
$ws.Cells.Item(15, 37).Value = 37
$ws.Cells.Item(15, 38).Value = 24
$ws.Cells.Item(15, 39).Value = 29
$ws.Cells.Item(15, 40).Value = 1.06
$ws.Cells.Item(15, 41).Value = 7.3

# Row 16
$ws.Cells.Item(16, 7).Value = 2.55
$ws.Cells.Item(16, 8).Value = 2.47
$ws.Cells.Item(16, 9).Value = 3.3
$ws.Cells.Item(16, 10).Value = 3.3
$ws.Cells.Item(16, 11).Value = 1.8
$ws.Cells.Item(16, 12).Value = 4
$ws.Cells.Item(16, 13).Value = 1.53
$ws.Cells.Item(16, 14).Value = 2.35
$ws.Cells.Item(16, 15).Value = 2.55
$ws.Cells.Item(16, 16).Value = 1.45
$ws.Cells.Item(16, 17).Value = 4.5
$ws.Cells.Item(16, 18).Value = 1.16
$ws.Cells.Item(16, 20).Value = 2.25
$ws.Cells.Item(16, 21).Value = 1.98
$ws.Cells.Item(16, 22).Value = 1.75
$ws.Cells.Item(16, 23).Value = 6.3
$ws.Cells.Item(16, 24).Value = 12
$ws.Cells.Item(16, 26).Value = 32
$ws.Cells.Item(16, 27).Value = 26
$ws.Cells.Item(16, 29).Value = 4.9
$ws.Cells.Item(16, 30).Value = 5
$ws.Cells.Item(16, 31).Value = 15
$ws.Cells.Item(16, 32).Value = 90
$ws.Cells.Item(16, 33).Value = 900
$ws.Cells.Item(16, 34).Value = 7.5
$ws.Cells.Item(16, 35).Value = 17
$ws.Cells.Item(16, 37).Value = 55
$ws.Cells.Item(16, 40).Value = 1.14
$ws.Cells.Item(16, 41).Value = 4.9
